$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.431881308555603
$ws.Range("B1").Value = 2.236846923828125
$ws.Range("C1").Value = 5.026980400085449
$ws.Range("D1").Value = 3.239110708236694
$ws.Range("E1").Value = 1.15222704410553
